$d = $word.ActiveDocument

# Locate the "Requisitos" heading paragraph, then the paragraph right
# after it, which holds the bulleted list of requirement lines. Each
# line lives in its own run, separated by manual line breaks (<w:br/>).
$paras = $d.Paragraphs
$headingIndex = -1
for ($i = 1; $i -le $paras.Count; $i++) {
    $t = $paras.Item($i).Range.Text.Trim()
    if ($t -eq "Requisitos") {
        $headingIndex = $i
    }
}

$listPara = $paras.Item($headingIndex + 1)
$rStart = $listPara.Range.Start
$fullText = $listPara.Range.Text

# Manual line breaks show up as vertical-tab (chr 11) inside Range.Text.
# Split the paragraph into the (start, end) offsets of each line, using
# those boundaries (the trailing paragraph mark is not a line).
$vtab = [char]11
$lineRanges = New-Object System.Collections.ArrayList
$lineStart = 0
for ($i = 0; $i -lt $fullText.Length; $i++) {
    $ch = $fullText.Substring($i, 1)
    if ($ch -eq $vtab) {
        $absStart = $rStart + $lineStart
        $absEnd = $rStart + $i
        $lineRanges.Add(@($absStart, $absEnd)) | Out-Null
        $lineStart = $i + 1
    }
}

# Find the line whose text is the "Física IV" requisite - this is the
# one that needs to move from its current spot to the front of the list.
$moveLineIndex = -1
for ($i = 0; $i -lt $lineRanges.Count; $i++) {
    $pair = $lineRanges[$i]
    $txt = $d.Range($pair[0], $pair[1]).Text
    if ($txt -like "LOB1021*") {
        $moveLineIndex = $i
    }
}

$movePair = $lineRanges[$moveLineIndex]
$moveText = $d.Range($movePair[0], $movePair[1]).Text

# Insert a brand new run (text + manual line break) at the very start of
# the list paragraph, using InsertBefore so it doesn't get merged into
# the neighbouring run even though formatting is identical.
$insertionPoint = $d.Range($rStart, $rStart)
$insertionPoint.InsertBefore($moveText + $vtab)

# Figure out the (shifted) position of the original line we just copied
# to the front, and remove it - including its own trailing line break -
# from its old location.
$shift = $moveText.Length + 1
$oldStart = $movePair[0] + $shift
$oldEnd = $movePair[1] + $shift + 1
$oldRange = $d.Range($oldStart, $oldEnd)
$oldRange.Delete()
